$d = $word.ActiveDocument

# Replace the representative's name and title in the "Narodni knihovna"
# party clause: "Ing. Petrem Knizkem, namestkem pro sekci digitalizace a
# technologie," becomes "Bc. Petrou Burdovou namestkyni pro sekci
# Digitalizace a technologie,"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Ing. Petrem Knížkem, náměstkem pro sekci digitalizace a technologie,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Bc. Petrou Burdovou náměstkyní pro sekci Digitalizace a technologie,",
    2)

Write-Host "Replace done"
